$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ "D" = "64.500.81"; "E" = "  -3.33%  " }
    3 = @{ "D" = "3.168.63"; "E" = "  -2.03%  " }
    4 = @{ "E" = "  +0.06%  " }
    5 = @{ "D" = "606.18"; "E" = "  +0.08%  " }
    6 = @{ "D" = "148.55"; "E" = "  -5.98%  " }
    7 = @{ "E" = "  +0.09%  " }
    8 = @{ "D" = "3.168.01"; "E" = "  -1.95%  " }
    9 = @{ "D" = "0.530"; "E" = "  -2.82%  " }
    10 = @{ "D" = "0.153"; "E" = "  -5.70%  " }
    11 = @{ "D" = "5.61"; "E" = "  -2.17%  " }
    12 = @{ "D" = "0.478"; "E" = "  -5.31%  " }
    13 = @{ "D" = "0.0000259"; "E" = "  -4.86%  " }
    14 = @{ "D" = "36.61"; "E" = "  -6.42%  " }
    15 = @{ "D" = "3.683.93"; "E" = "  -2.16%  " }
    16 = @{ "D" = "64.500.06"; "E" = "  -3.31%  " }
    17 = @{ "E" = "  +1.19%  " }
    18 = @{ "D" = "3.165.11"; "E" = "  -2.12%  " }
    19 = @{ "D" = "6.98"; "E" = "  -4.34%  " }
    20 = @{ "D" = "483.18"; "E" = "  -5.35%  " }
    21 = @{ "D" = "14.64"; "E" = "  -4.48%  " }
    22 = @{ "D" = "0.714"; "E" = "  -3.38%  " }
    23 = @{ "D" = "7.81"; "E" = "  -3.31%  " }
    24 = @{ "D" = "13.84"; "E" = "  -5.66%  " }
    25 = @{ "D" = "83.69"; "E" = "  -2.56%  " }
    26 = @{ "E" = "  -0.04%  " }
    27 = @{ "D" = "2.91"; "E" = "  -3.35%  " }
    28 = @{ "D" = "8.57"; "E" = "  -6.17%  " }
    29 = @{ "D" = "2.23"; "E" = "  -6.05%  " }
    30 = @{ "E" = "  -20.16%  " }
    31 = @{ "D" = "6.97"; "E" = "  -1.20%  " }
    32 = @{ "D" = "2.78"; "E" = "  -5.08%  " }
    33 = @{ "E" = "  -0.01%  " }
    34 = @{ "D" = "26.46"; "E" = "  -6.80%  " }
    35 = @{ "E" = "  -5.72%  " }
    36 = @{ "D" = "6.10"; "E" = "  -5.75%  " }
    37 = @{ "D" = "54.43"; "E" = "  -1.95%  " }
    38 = @{ "D" = "0.0₃0741"; "E" = "  -8.68%  " }
    39 = @{ "E" = "  -2.08%  " }
    40 = @{ "D" = "458.55"; "E" = "  -8.82%  " }
    41 = @{ "B" = "VeChain"; "C" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; "D" = "0.0400"; "E" = "  -5.33%  " }
    42 = @{ "B" = "Kaspa"; "C" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; "D" = "0.124"; "E" = "  -3.44%  " }
    43 = @{ "D" = "8.49"; "E" = "  -2.92%  " }
    44 = @{ "D" = "2.895.00"; "E" = "  -1.98%  " }
    45 = @{ "D" = "0.273"; "E" = "  -8.18%  " }
    46 = @{ "D" = "2.30"; "E" = "  -7.31%  " }
    47 = @{ "D" = "26.63"; "E" = "  -6.48%  " }
    49 = @{ "D" = "2.35"; "E" = "  -4.32%  " }
    50 = @{ "D" = "0.115"; "E" = "  -2.43%  " }
    51 = @{ "D" = "118.61"; "E" = "  -2.62%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
